$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we are about to update so that
# numeric-looking / percent-looking strings are preserved as text
# (matching the original inlineStr cell type), instead of being
# auto-converted into numbers/percentages by Excel.
$updatedCells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","E20","E21","D22","E22","E23","E24","D25","E25","E26","D39","E39","D40","E40","D41","E41","D42","E42","D43","D44","E44","D45","E45","D46","E46","E47","D48","E48","E49","E50","E51")
foreach ($addr in $updatedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the new values exactly as strings.
$ws.Range("D2").Value = "288.25"
$ws.Range("E2").Value = "-0.80%"
$ws.Range("D3").Value = "30.98"
$ws.Range("E3").Value = "0.78%"
$ws.Range("D4").Value = "4.924"
$ws.Range("E4").Value = "-0.58%"
$ws.Range("D5").Value = "0.07373"
$ws.Range("E5").Value = "2.23%"
$ws.Range("D6").Value = "2.256"
$ws.Range("E6").Value = "25.38%"
$ws.Range("D7").Value = "7.701"
$ws.Range("E7").Value = "0.40%"
$ws.Range("D8").Value = "3.733"
$ws.Range("E8").Value = "0.25%"
$ws.Range("D9").Value = "0.9077"
$ws.Range("E9").Value = "1.28%"
$ws.Range("D10").Value = "0.08751"
$ws.Range("E10").Value = "13.53%"
$ws.Range("D11").Value = "0.1685"
$ws.Range("E11").Value = "1.92%"
$ws.Range("D12").Value = "0.08208"
$ws.Range("E12").Value = "2.90%"
$ws.Range("D13").Value = "0.03114"
$ws.Range("E13").Value = "2.11%"
$ws.Range("D14").Value = "0.09948"
$ws.Range("E14").Value = "-0.74%"
$ws.Range("D15").Value = "0.001497"
$ws.Range("E15").Value = "0.20%"
$ws.Range("D16").Value = "0.005716"
$ws.Range("E16").Value = "0.99%"
$ws.Range("D17").Value = "3.491"
$ws.Range("E17").Value = "0.43%"
$ws.Range("D18").Value = "2.082"
$ws.Range("E18").Value = "-0.10%"
$ws.Range("E20").Value = "-1.37%"
$ws.Range("E21").Value = "-5.10%"
$ws.Range("D22").Value = "0.2123"
$ws.Range("E22").Value = "1.11%"
$ws.Range("E23").Value = "0.86%"
$ws.Range("E24").Value = "-0.21%"
$ws.Range("D25").Value = "0.004141"
$ws.Range("E25").Value = "3.06%"
$ws.Range("E26").Value = "4.10%"
$ws.Range("D39").Value = "0.01582"
$ws.Range("E39").Value = "-1.03%"
$ws.Range("D40").Value = "0.04465"
$ws.Range("E40").Value = "1.68%"
$ws.Range("D41").Value = "0.007317"
$ws.Range("E41").Value = "-0.21%"
$ws.Range("D42").Value = "0.009562"
$ws.Range("E42").Value = "24.24%"
$ws.Range("D43").Value = "0.1324"
$ws.Range("D44").Value = "0.002252"
$ws.Range("E44").Value = "9.33%"
$ws.Range("D45").Value = "0.008434"
$ws.Range("E45").Value = "-8.48%"
$ws.Range("D46").Value = "0.00006136"
$ws.Range("E46").Value = "4.98%"
$ws.Range("E47").Value = "0.11%"
$ws.Range("D48").Value = "2.190"
$ws.Range("E48").Value = "-2.48%"
$ws.Range("E49").Value = "-33.26%"
$ws.Range("E50").Value = "0.11%"
$ws.Range("E51").Value = "0.11%"
